{"js": "// The Fig1B \"input:\" paragraph cites a filename that needs updating:\n//   allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_biotype_06242025.csv\n// becomes\n//   allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_genetype_06302025.csv\n// (i.e. \"bio\" -> \"gene\" and \"24\" -> \"30\"), and - matching how Word itself\n// lays things out after a couple of in-place retypes - the new text ends\n// up split across five highlighted runs instead of the original single run.\n\nconst oldText = \"allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_biotype_06242025.csv\";\n\nconst body = context.document.body;\nconst hits = body.search(oldText, { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find target filename text to update.\");\n}\n\n// Flat-OPC package wrapping just the fragment of document.xml we want to\n// splice in, so the single highlighted run becomes the five runs Word\n// produced: \"...modu_\" + \"gene\" + \"type_06\" + \"30\" + \"2025.csv\".\nconst replacementOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r>\n              <w:rPr><w:highlight w:val=\"yellow\"/></w:rPr>\n              <w:t>allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:highlight w:val=\"yellow\"/></w:rPr>\n              <w:t>gene</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:highlight w:val=\"yellow\"/></w:rPr>\n              <w:t>type_06</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:highlight w:val=\"yellow\"/></w:rPr>\n              <w:t>30</w:t>\n            </w:r>\n            <w:r>\n              <w:rPr><w:highlight w:val=\"yellow\"/></w:rPr>\n              <w:t>2025.csv</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n// Replace the whole matched range's contents with the 5-run fragment above.\nhits.items[0].insertOoxml(replacementOoxml, \"Replace\");\nawait context.sync();\n", "ps1": "# The Fig1B \"input:\" paragraph cites a filename that needs updating:\n#   allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_biotype_06242025.csv\n# becomes\n#   allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_genetype_06302025.csv\n# (i.e. \"bio\" -> \"gene\" and \"24\" -> \"30\"), and - matching how Word itself\n# lays things out after a couple of in-place retypes - the new text ends\n# up split across five highlighted runs instead of the original single run.\n\n$d = $word.ActiveDocument\n\n$oldText = \"allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_biotype_06242025.csv\"\n\n$hit = $d.Content\n$find = $hit.Find\n$find.Text = $oldText\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw \"Could not find target filename text to update.\"\n}\n\n# $hit now spans exactly the matched text (Find collapses/extends its own\n# range in place). Re-seat onto a plain Range over that same span (rather\n# than reusing the Find-bound range object) so InsertXML REPLACES this\n# range's contents instead of inserting after it.\n$r = $d.Range($hit.Start, $hit.End)\n\n# Flat-OPC package wrapping just the fragment of document.xml we want to\n# splice in, so the single highlighted run becomes the five runs Word\n# produced: \"...modu_\" + \"gene\" + \"type_06\" + \"30\" + \"2025.csv\".\n$replacementXml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n'<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n      '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' +\n          '<w:p>' +\n            '<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>allgenes_XAKRIPr_XAKnetworkr_common_sil_modu_</w:t></w:r>' +\n            '<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>gene</w:t></w:r>' +\n            '<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>type_06</w:t></w:r>' +\n            '<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>30</w:t></w:r>' +\n            '<w:r><w:rPr><w:highlight w:val=\"yellow\"/></w:rPr><w:t>2025.csv</w:t></w:r>' +\n          '</w:p>' +\n        '</w:body>' +\n      '</w:document>' +\n    '</pkg:xmlData>' +\n  '</pkg:part>' +\n'</pkg:package>'\n\n$r.InsertXML($replacementXml)\n"}
